# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Update the "Good Roaming Calculation (%)" value for the first bad driver row
$ws.Range("D3").Value = 82.09999999999999

# Fill in the "Driver Vintage" date for the previously blank row.
# A leading apostrophe forces literal text so Excel doesn't auto-convert the
# ISO-looking string into a date serial; then re-apply the original
# (General, right-aligned) number format from a sibling cell so the cell's
# style index is unchanged.
$ws.Range("E12").Value = "'2022-08-29"
$ws.Range("D12").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
